# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - "Ready for handoff" -> "Handed back: in sync with en-US" on all three sheets
#  - Populates "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
#    for both rows on the zh-cn and de-de sheets (with hyperlinks on the target-file column)
#  - Widens a few columns that now hold longer content

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdA = "34dc4ef2-d55c-4e4a-8515-ffb9baa38d99.md"
$mdB = "c03eaaa9-6ef0-4fde-8007-3eed47ecb17f.md"
$urlA = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/421977e98050aa6e4e4c6e0ac3733669b0c4cfd8/e2e/34dc4ef2-d55c-4e4a-8515-ffb9baa38d99.md"
$urlB = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/421977e98050aa6e4e4c6e0ac3733669b0c4cfd8/e2e/c03eaaa9-6ef0-4fde-8007-3eed47ecb17f.md"

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $newStatus
$wsOverview.Range("F2").Value2 = $newStatus
$wsOverview.Range("E3").Value2 = $newStatus
$wsOverview.Range("F3").Value2 = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value2 = $newStatus
$wsZh.Range("C3").Value2 = $newStatus

$wsZh.Range("I2").Value2 = $mdA
$wsZh.Range("J2").Value2 = "34dc4ef2-d55c-4e4a-8515-ffb9baa38d99.09ee360233700a4462c9db4da6df630a27d552de.zh-cn.xlf"
$wsZh.Range("K2").Value2 = "2016-09-07 01:05:32"

$wsZh.Range("I3").Value2 = $mdB
$wsZh.Range("J3").Value2 = "c03eaaa9-6ef0-4fde-8007-3eed47ecb17f.2ed4310148090bf4924941f61a08b9b27ea0c269.zh-cn.xlf"
$wsZh.Range("K3").Value2 = "2016-09-07 01:05:32"

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $urlA, "", "", $mdA)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $urlB, "", "", $mdB)

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value2 = $newStatus
$wsDe.Range("C3").Value2 = $newStatus

$wsDe.Range("I2").Value2 = $mdA
$wsDe.Range("J2").Value2 = "34dc4ef2-d55c-4e4a-8515-ffb9baa38d99.09ee360233700a4462c9db4da6df630a27d552de.de-de.xlf"
$wsDe.Range("K2").Value2 = "2016-09-07 01:05:40"

$wsDe.Range("I3").Value2 = $mdB
$wsDe.Range("J3").Value2 = "c03eaaa9-6ef0-4fde-8007-3eed47ecb17f.2ed4310148090bf4924941f61a08b9b27ea0c269.de-de.xlf"
$wsDe.Range("K3").Value2 = "2016-09-07 01:05:40"

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $urlA, "", "", $mdA)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $urlB, "", "", $mdB)

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

Write-Host "Handback report generated."
